# Add season record columns (Wins, Losses, Ties) to the roster/stats sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Create the three new header cells (AD1:AF1) using the same visual style
# as the existing header row (copy formatting from an existing header cell,
# then overwrite the text) so the new headers look consistent with "WAR",
# "Salary", etc.
$ws.Range("AA1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record for every player row (2-49) with the team's
# season record: 68 wins, 94 losses, 0 ties.
for ($r = 2; $r -le 49; $r++) {
  $ws.Cells.Item($r, 30).Value = 68
  $ws.Cells.Item($r, 31).Value = 94
  $ws.Cells.Item($r, 32).Value = 0
}

Write-Output "Added Wins/Losses/Ties columns (AD:AF) for rows 1-49"
